$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 had no data before (sheet dimension started at row 2); just populate
# it directly instead of inserting/shifting existing rows.

# Row height for the new header row.
$ws.Rows.Item(1).RowHeight = 25

# Header values.
$ws.Range("A1").Value = "Description"
$ws.Range("B1").Value = "Qty"
$ws.Range("C1").Value = "Net Weight"
$ws.Range("D1").Value = "Gross Weight"

# Formatting to match the other label/value columns.
$ws.Range("A1").HorizontalAlignment = -4131
$ws.Range("A1").VerticalAlignment = -4108

$ws.Range("B1").HorizontalAlignment = -4108
$ws.Range("B1").VerticalAlignment = -4108

$ws.Range("C1:D1").HorizontalAlignment = -4152
$ws.Range("C1:D1").VerticalAlignment = -4108
